$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1916468203859388
$ws.Range("C2").Value = 0.7358528560593535
$ws.Range("D2").Value = 0.8902295501318038
$ws.Range("E2").Value = 0.9435197666884376
$ws.Range("F2").Value = 0.9587257026967815
$ws.Range("G2").Value = 14
